# New weekly price record for Jengibre (Terminal La Palmera de La Serena)
# is inserted as a new data row right after the existing row 74, pushing
# every subsequent row down by one (old row 75 -> new row 76, ..., old
# row 181 -> new row 182).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; everything from 75 downward shifts down.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly record.
$ws.Cells.Item(75, 1).Value = 8
$ws.Cells.Item(75, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).Value = 45175
$ws.Cells.Item(75, 5).Value = 4
$ws.Cells.Item(75, 6).Value = 100114007
$ws.Cells.Item(75, 7).Value = "Jengibre"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 480
$ws.Cells.Item(75, 11).Value = 17000
$ws.Cells.Item(75, 12).Value = 18000
$ws.Cells.Item(75, 13).Value = 17500
$ws.Cells.Item(75, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(75, 15).Value = "Perú"
$ws.Cells.Item(75, 16).Value = 1346
$ws.Cells.Item(75, 17).Value = 13
$ws.Cells.Item(75, 18).Value = "Hortaliza"
